$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AA2" = 29
    "AB2" = 35
    "AC2" = 8.75
    "AD2" = 5.6
    "AE2" = 12.5
    "AF2" = 55
    "AH2" = 5.5
    "AI2" = 6.9
    "AJ2" = 7.1
    "AK2" = 12.5
    "AL2" = 12.5
    "AM2" = 23
    "G2" = 3.9
    "H2" = 3.25
    "I2" = 1.83
    "J2" = 4.35
    "K2" = 2.02
    "L2" = 2.47
    "M2" = 1.36
    "N2" = 2.9
    "O2" = 1.98
    "P2" = 1.65
    "Q2" = 3.7
    "R2" = 1.24
    "S2" = 1.42
    "T2" = 2.37
    "U2" = 1.97
    "W2" = 9
    "X2" = 17.5
    "Y2" = 11
    "Z2" = 45
    "M3" = 1.33
    "N3" = 3.05
    "Q3" = 3.5
    "R3" = 1.26
    "S3" = 1.44
    "T3" = 2.6
    "U3" = 1.82
    "V3" = 1.89
    "AH4" = 6.7
    "AI4" = 10.25
    "AL4" = 17.5
    "AC5" = 12
    "AD5" = 6.5
    "AH5" = 9
    "H5" = 3.4
    "O5" = 1.75
    "P5" = 2.05
    "S5" = 1.36
    "T5" = 3
    "AA6" = 18
    "AC6" = 10.25
    "AD6" = 6.4
    "AH6" = 9.25
    "AI6" = 15
    "AK6" = 35
    "G6" = 2.3
    "I6" = 2.87
    "J6" = 2.87
    "K6" = 2.07
    "L6" = 3.5
    "W6" = 8.5
    "X6" = 12
    "Y6" = 9
    "Z6" = 24
    "AA7" = 24
    "AD7" = 6.8
    "AG7" = 450
    "AH7" = 8.25
    "AI7" = 10.75
    "AK7" = 20
    "AL7" = 17
    "G7" = 2.95
    "H7" = 3.5
    "I7" = 2.15
    "J7" = 3.45
    "L7" = 2.72
    "O7" = 1.75
    "P7" = 1.85
    "W7" = 10
    "X7" = 15.5
    "Y7" = 10.75
    "Z7" = 35
    "M8" = 1.29
    "N8" = 3.25
    "Q8" = 3.2
    "R8" = 1.3
    "U8" = 1.76
    "V8" = 1.96
    "AC11" = 15
    "AE11" = 12
    "AI11" = 12
    "AP11" = 2.03
    "AQ11" = 1.78
    "G11" = 2.9
    "I11" = 2.15
    "J11" = 3.4
    "L11" = 2.75
    "W11" = 13
    "X11" = 17
    "AB13" = 26
    "AC13" = 11
    "AJ13" = 12
    "G13" = 2.2
    "H13" = 3.3
    "I13" = 3
    "J13" = 2.88
    "L13" = 3.5
    "O13" = 1.9
    "P13" = 1.9
    "X13" = 11
    "Y13" = 9.5
    "Z13" = 21
    "AC14" = 6.5
    "AI14" = 26
    "AJ14" = 19
    "AM14" = 51
    "AR14" = 1.85
    "AS14" = 1.95
    "G14" = 1.67
    "I14" = 5.5
    "J14" = 2.4
    "Q14" = 5
    "R14" = 1.17
    "Y14" = 9.5
    "AH16" = 12
    "AJ16" = 13
    "I16" = 4.1
    "L16" = 4.5
    "X16" = 9
    "AC18" = 8.5
    "AE18" = 30
    "AH18" = 16
    "AI18" = 55
    "AJ18" = 29
    "AK18" = 250
    "G18" = 1.36
    "H18" = 4.25
    "I18" = 8.25
    "J18" = 1.88
    "L18" = 7.7
    "M18" = 1.33
    "N18" = 2.77
    "O18" = 1.98
    "P18" = 1.65
    "Q18" = 3.25
    "R18" = 1.24
    "T18" = 2.45
    "U18" = 2.4
    "W18" = 5.1
    "X18" = 5.3
    "Z18" = 7.8
    "AA19" = 30
    "AC19" = 6.1
    "AE19" = 18
    "AF19" = 120
    "AH19" = 6.4
    "AI19" = 11.75
    "AJ19" = 10.75
    "AK19" = 32
    "AL19" = 29
    "AM19" = 50
    "AN19" = 1.12
    "G19" = 2.75
    "H19" = 2.85
    "I19" = 2.67
    "J19" = 3.55
    "K19" = 1.82
    "L19" = 3.45
    "M19" = 1.52
    "N19" = 2.22
    "O19" = 2.47
    "P19" = 1.42
    "Q19" = 4.25
    "R19" = 1.14
    "T19" = 2.1
    "U19" = 2.05
    "V19" = 1.62
    "W19" = 6.5
    "X19" = 12
    "Y19" = 11
    "Z19" = 32
    "AL20" = 13
    "K20" = 2.12
    "U20" = 2.07
    "V20" = 1.6
    "AB21" = 37
    "AC21" = 10.5
    "AE21" = 28
    "AF21" = 175
    "AH21" = 19
    "AI21" = 60
    "AJ21" = 29
    "AK21" = 250
    "AL21" = 120
    "AM21" = 120
    "G21" = 1.32
    "H21" = 4.65
    "I21" = 8.75
    "K21" = 2.3
    "L21" = 7.7
    "M21" = 1.26
    "N21" = 3.15
    "O21" = 1.78
    "P21" = 1.82
    "Q21" = 2.82
    "R21" = 1.32
    "S21" = 1.38
    "T21" = 2.6
    "U21" = 2.25
    "V21" = 1.5
    "W21" = 5.8
    "X21" = 5.5
    "Y21" = 9
    "Z21" = 7.5
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
